# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. This updates the DAMSLTag (column I) and
# DialogAct (column J) values for a set of rows in the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 6;   I = "sv"; J = "Statement-opinion" },
    @{ Row = 8;   I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 15;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 18;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 24;  I = "%";  J = "Uninterpretable" },
    @{ Row = 64;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 75;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 82;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 83;  I = "%";  J = "Uninterpretable" },
    @{ Row = 85;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 88;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 94;  I = "%";  J = "Uninterpretable" },
    @{ Row = 95;  I = "%";  J = "Uninterpretable" },
    @{ Row = 103; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 104; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 109; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 112; I = "sd"; J = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
